$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.811.71'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '2.581.94'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'582.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = "'144.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = "'0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = "'27.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = '3.047.35'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '62.744.13'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '2.588.48'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").Value = "'11.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'339.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = "'4.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = "'6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = "'5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.17%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'67.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.61%  '
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = "'8.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("B28").Value = 'SuiNetwork'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D28").Value = "'1.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = "'8.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = "'1.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").Value = '0.0₃0809'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = "'453.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.80%  '
$ws.Range("D34").Value = "'176.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("D37").Value = "'0.400"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").Value = "'4.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("D42").Value = "'157.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.49%  '
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("D44").Value = "'21.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.44%  '
$ws.Range("D45").Value = "'0.630"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.74%  '
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").Value = "'0.0965"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").Value = "'0.0234"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("D49").Value = "'18.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").Value = "'11.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("E51").Value = '  -1.04%  '
